# Auto-generated Excel COM-interop edit script
# Applies numeric value updates (and a few cell clears) to the
# Brynhildr_Profits workbook's leve-profit computation columns
# (H..N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 2427
$ws.Range("J58").Value = 2549.625
$ws.Range("L58").Value = 7648.875
$ws.Range("N58").Value = -7948.875
$ws.Range("H96").Value = 2817.182
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 2817.182
$ws.Range("K96").Value = 0
$ws.Range("L96").ClearContents()
$ws.Range("M96").Value = 8451.545999999998
$ws.Range("N96").Value = -11197.546
$ws.Range("H112").Value = 2509
$ws.Range("J112").Value = 3000
$ws.Range("L112").Value = 9000
$ws.Range("N112").Value = -11216

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1430.28
$ws.Range("I32").Value = 1218.6702
$ws.Range("J32").Value = 4745.5
$ws.Range("K32").Value = 1218.6702
$ws.Range("L32").Value = 4745.5
$ws.Range("M32").Value = -931.6702
$ws.Range("N32").Value = -5319.5
$ws.Range("H61").Value = 2373.3333
$ws.Range("I61").Value = 2043.25
$ws.Range("K61").Value = 2043.25
$ws.Range("M61").Value = -1831.25
$ws.Range("H74").Value = 8928.25
$ws.Range("J74").Value = 14206.071
$ws.Range("L74").Value = 14206.071
$ws.Range("N74").Value = -15954.071
$ws.Range("H77").Value = 8928.25
$ws.Range("J77").Value = 14206.071
$ws.Range("L77").Value = 71030.355
$ws.Range("N77").Value = -79766.355
$ws.Range("H132").Value = 3016.6365
$ws.Range("I132").Value = 2164.7856
$ws.Range("K132").Value = 6494.3568
$ws.Range("M132").Value = -3964.3568
$ws.Range("H136").Value = 2373.3333
$ws.Range("I136").Value = 2043.25
$ws.Range("K136").Value = 6129.75
$ws.Range("M136").Value = -3579.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 8960.714
$ws.Range("I99").Value = 9419.308000000001
$ws.Range("K99").Value = 9419.308000000001
$ws.Range("M99").Value = -7921.308000000001
$ws.Range("H105").Value = 2139.342
$ws.Range("I105").Value = 1845.8695
$ws.Range("J105").Value = 2589.3333
$ws.Range("K105").Value = 1845.8695
$ws.Range("L105").Value = 2589.3333
$ws.Range("M105").Value = -98.86950000000002
$ws.Range("N105").Value = -6083.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 728.4074000000001
$ws.Range("I107").Value = 574.9524
$ws.Range("K107").Value = 574.9524
$ws.Range("M107").Value = 1345.0476
$ws.Range("H115").Value = 12290
$ws.Range("J115").Value = 12290
$ws.Range("L115").Value = 12290
$ws.Range("N115").Value = -14640
$ws.Range("H121").Value = 58625
$ws.Range("J121").Value = 58625
$ws.Range("L121").Value = 58625
$ws.Range("N121").Value = -61245
$ws.Range("H132").Value = 4587.24
$ws.Range("I132").Value = 4245.6904
$ws.Range("K132").Value = 12737.0712
$ws.Range("M132").Value = -10207.0712

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 126.44444
$ws.Range("I2").Value = 149
$ws.Range("J2").Value = 115.166664
$ws.Range("K2").Value = 894
$ws.Range("L2").Value = 690.999984
$ws.Range("M2").Value = -781
$ws.Range("N2").Value = -916.999984
$ws.Range("H86").Value = 1900
$ws.Range("I86").Value = 1900
$ws.Range("K86").Value = 5700
$ws.Range("M86").Value = -4514
$ws.Range("H89").Value = 1900
$ws.Range("I89").Value = 1900
$ws.Range("K89").Value = 17100
$ws.Range("M89").Value = -11172
$ws.Range("H127").Value = 4966.5
$ws.Range("J127").Value = 4966.5
$ws.Range("L127").Value = 14899.5
$ws.Range("N127").Value = -24819.5
$ws.Range("H134").Value = 2954.6667
$ws.Range("I134").Value = 2091.0476
$ws.Range("K134").Value = 6273.1428
$ws.Range("M134").Value = -1203.1428

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").ClearContents()
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = 0
$ws.Range("H35").Value = 13166.667
$ws.Range("I35").Value = 11250
$ws.Range("J35").Value = 17000
$ws.Range("K35").Value = 11250
$ws.Range("L35").Value = 17000
$ws.Range("M35").Value = -10952
$ws.Range("N35").Value = -17596
$ws.Range("H70").Value = 16315.412
$ws.Range("I70").Value = 17486.357
$ws.Range("K70").Value = 17486.357
$ws.Range("M70").Value = -17216.357
$ws.Range("H73").Value = 16315.412
$ws.Range("I73").Value = 17486.357
$ws.Range("K73").Value = 17486.357
$ws.Range("M73").Value = -16550.357
$ws.Range("H132").Value = 12408.526
$ws.Range("I132").Value = 14240.9375
$ws.Range("K132").Value = 42722.8125
$ws.Range("M132").Value = -40192.8125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 8246.267
$ws.Range("I61").Value = 9207.916999999999
$ws.Range("J61").Value = 4399.6665
$ws.Range("K61").Value = 9207.916999999999
$ws.Range("L61").Value = 4399.6665
$ws.Range("M61").Value = -9005.916999999999
$ws.Range("N61").Value = -4803.6665
$ws.Range("H68").Value = 12725.179
$ws.Range("I68").Value = 10183.277
$ws.Range("K68").Value = 10183.277
$ws.Range("M68").Value = -9434.277
$ws.Range("H71").Value = 12725.179
$ws.Range("I71").Value = 10183.277
$ws.Range("K71").Value = 50916.385
$ws.Range("M71").Value = -47172.385
$ws.Range("H93").Value = 4890.6
$ws.Range("I93").Value = 2536.1667
$ws.Range("J93").Value = 8422.25
$ws.Range("K93").Value = 2536.1667
$ws.Range("L93").Value = 8422.25
$ws.Range("M93").Value = -1288.1667
$ws.Range("N93").Value = -10918.25
$ws.Range("H113").Value = 8246.267
$ws.Range("I113").Value = 9207.916999999999
$ws.Range("J113").Value = 4399.6665
$ws.Range("K113").Value = 9207.916999999999
$ws.Range("L113").Value = 4399.6665
$ws.Range("M113").Value = -7037.916999999999
$ws.Range("N113").Value = -8739.666499999999
$ws.Range("H132").Value = 1526.6
$ws.Range("I132").Value = 1526.6
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4579.799999999999
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -2049.799999999999
$ws.Range("H136").Value = 4248.9375
$ws.Range("I136").Value = 4077.2307
$ws.Range("J136").Value = 4993
$ws.Range("K136").Value = 12231.6921
$ws.Range("L136").Value = 14979
$ws.Range("M136").Value = -9681.6921
$ws.Range("N136").Value = -20079
$ws.Range("H140").Value = 94925
$ws.Range("J140").Value = 94925
$ws.Range("L140").Value = 94925
$ws.Range("N140").Value = -105285

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 578570.3
$ws.Range("I4").Value = 1999996.5
$ws.Range("K4").Value = 1999996.5
$ws.Range("M4").Value = -1999883.5
$ws.Range("H132").Value = 2043.25
$ws.Range("I132").Value = 1649.1052
$ws.Range("K132").Value = 4947.3156
$ws.Range("M132").Value = -2417.3156
$ws.Range("H136").Value = 1830.12
$ws.Range("J136").Value = 1544.091
$ws.Range("L136").Value = 4632.272999999999
$ws.Range("N136").Value = -9732.272999999999
